$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the order of Alex Francoletti / Kevin, and expand "Kevin" to "Kevin Chen"
$ws.Range("B16").Value = "Kevin Chen"
$ws.Range("B17").Value = "Alex Francoletti"

# Flesh out the Music section (previously just "Sources" / "Music" headers)
$ws.Range("A22").Value = "Music (Selection)"
$ws.Range("B22").Value = "James Laks"

$ws.Range("A23").Value = "Music (Sources)"
$ws.Range("B23").Value = "Bensound"

$ws.Range("A24").ClearContents()
$ws.Range("B24").Value = "Purple Planet Music"

$ws.Range("A25").ClearContents()
$ws.Range("B25").ClearContents()

# Move selection to reflect where editing left off
$ws.Range("C22").Select()
